# Quarterly indexing esoteric bug-fix operation:
# Column A (rows 2..N) holds the "as-of" date for each forecast row. Each
# date had been stamped as the 1st of its month; the fix re-stamps it to
# the 15th of the *following* month (the corrected quarterly indexing
# convention used elsewhere in the pipeline).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-ShiftedSerial($serial) {
    $d = [DateTime]::FromOADate($serial)
    $year = $d.Year
    $month = $d.Month + 1
    if ($month -gt 12) {
        $month = 1
        $year = $year + 1
    }
    $newDate = Get-Date -Year $year -Month $month -Day 15 -Hour 0 -Minute 0 -Second 0
    return [Math]::Floor($newDate.ToOADate())
}

# Find the last used row in column A (xlUp = -4162), falling back to 63
# (the known extent of this sheet) if that lookup fails for any reason.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 63 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value2 = Get-ShiftedSerial $old
    }
}
